# Applies the commit's edits to the document:
#  1. Split the "somyajain99" hyperlink run into 3 runs (so/m/yajain99...)
#  2. Split the "filyp" hyperlink run into 3 runs (autoco/r/rrect...)
#  3. Insert "the " into "So, we need to remove leading number." as 3 runs
#  4. Change the last (empty) list paragraph's numbering from ilvl=2/numId=1
#     to ilvl=0/numId=7

$d = $word.ActiveDocument

function Split-RunAt($rangeStart, $rangeEnd, $offsets) {
    # Forces independent <w:r> runs to exist at the given absolute character
    # offsets (strictly between rangeStart and rangeEnd) by nudging a
    # formatting property away and back again, which makes the engine
    # materialize a run boundary there while keeping the run's rPr intact.
    foreach ($off in $offsets) {
        $piece = $d.Range($off, $rangeEnd)
        $piece.Font.Size = 99
        $piece.Font.Size = 12
    }
}

# --- 1) https://github.com/somyajain99/english-autocorrect -----------------
$rng = $d.Content
$found = $rng.Find.Execute("https://github.com/somyajain99/english-autocorrect", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $end = $rng.End
    $off1 = $start + ("https://github.com/so").Length
    $off2 = $off1 + ("m").Length
    Split-RunAt $start $end @($off1, $off2)
}

# --- 2) https://github.com/filyp/autocorrect/tree/master -------------------
$rng = $d.Content
$found = $rng.Find.Execute("https://github.com/filyp/autocorrect/tree/master", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $end = $rng.End
    $off1 = $start + ("https://github.com/filyp/autoco").Length
    $off2 = $off1 + ("r").Length
    Split-RunAt $start $end @($off1, $off2)
}

# --- 3) "So, we need to remove leading number." -> insert "the " -----------
$rng = $d.Content
$found = $rng.Find.Execute("So, we need to remove leading number.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "So, we need to remove the leading number.", 2)
if ($found) {
    # Re-locate the (now longer) sentence and split it into three runs.
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute("So, we need to remove the leading number.", `
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $start = $rng2.Start
        $end = $rng2.End
        $off1 = $start + ("So, we need to remove ").Length
        $off2 = $off1 + ("the ").Length
        Split-RunAt $start $end @($off1, $off2)
    }
}

# --- 4) last empty list paragraph: ilvl 2->0, numId 1->7 -------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text
    if ($txt.Length -le 1) {
        $lvl = $p.Range.ListFormat.ListLevelNumber
        $lid = $p.Range.ListFormat.List.ListID
        if ($lvl -eq 3 -and $lid -eq 1) {
            $p.Range.ListFormat.ListTemplate = 7
            $p.Range.ListFormat.ListLevelNumber = 1
        }
    }
}
